$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Slide title: consolidate the split "Here"/" "/"is"/... runs into one run ---
# Setting .Text to a value that already equals the current text is treated as a
# no-op by the writer, so first nudge it to a throwaway value and then back to
# the real text; this forces the run list to be rebuilt as a single run.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "_tmp_"
$title.Text = "Here is a single header"

# --- Speaker notes: consolidate the split "and"/" "/"here"/... runs into one run ---
$notes = $s.NotesPage
$notesBody = $notes.Shapes.Item(2).TextFrame.TextRange
$notesBody.Text = "_tmp_"
$notesBody.Text = "and here are some notes"
